# Nexial "remote-amazon" scenario template update:
#  - add new JSON helper function storeKeys(json,jsonpath,var) to the
#    #system sheet's `json` lookup column (inserted alphabetically
#    between storeCount and storeValue)
#  - remove the now-obsolete `text` target category (a single-entry
#    list) which frees up its column, shifting the web/webalert/
#    webcookie/ws/ws.async/xml lookup columns one column to the left
#  - keep all named ranges in sync with their (possibly resized /
#    shifted) column ranges on the #system sheet
#
# NOTE: Range.Insert()/Range.Delete() on a single cell shift the whole
# row in this environment, not just the target column, so the in-
# column row shifts below are done by copying cell values manually
# (column-by-column) instead of relying on Insert/Delete for those.
# A full-column Delete (Columns(...).Delete()) does behave correctly
# (column-only), so that one is used as-is for dropping column Y.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("#system")

# 1) Insert the new JSON function in column M (the `json` list) at
#    row 16, pushing storeValue/storeValues down by one row. Column M
#    data currently runs from row 2 through row 17.
for ($r = 17; $r -ge 16; $r--) {
    $ws.Cells.Item($r + 1, 13).Value = $ws.Cells.Item($r, 13).Value()
}
$ws.Cells.Item(16, 13).Value = "storeKeys(json,jsonpath,var)"

# 2) Drop the `text` entry from column A (the `target` list) -- it
#    currently sits at A25 -- shifting everything below it up one row
#    (column A data runs through row 31).
for ($r = 25; $r -le 30; $r++) {
    $ws.Cells.Item($r, 1).Value = $ws.Cells.Item($r + 1, 1).Value()
}
$ws.Cells.Item(31, 1).Value = ""

# 3) Drop the now-unused `text` lookup column (Y) entirely, shifting
#    web / webalert / webcookie / ws / ws.async / xml one column left
#    (Z->Y, AA->Z, AB->AA, AC->AB, AD->AC, AE->AD). A full-column
#    delete correctly shifts only the affected columns.
$ws.Columns("Y:Y").Delete()

# 4) Re-point the named ranges that moved / resized as a result.
$wb.Names.Item("json").RefersTo = "='#system'!`$M`$2:`$M`$18"
$wb.Names.Item("target").RefersTo = "='#system'!`$A`$2:`$A`$30"
$wb.Names.Item("web").RefersTo = "='#system'!`$Y`$2:`$Y`$129"
$wb.Names.Item("webalert").RefersTo = "='#system'!`$Z`$2:`$Z`$8"
$wb.Names.Item("webcookie").RefersTo = "='#system'!`$AA`$2:`$AA`$8"
$wb.Names.Item("ws").RefersTo = "='#system'!`$AB`$2:`$AB`$17"
$wb.Names.Item("ws.async").RefersTo = "='#system'!`$AC`$2:`$AC`$8"
$wb.Names.Item("xml").RefersTo = "='#system'!`$AD`$2:`$AD`$27"
